$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-03 Sunday" "2024-11-04 Monday"
Replace-Text "72×39=" "64×41="
Replace-Text "82×34=" "60×42="
Replace-Text "94×81=" "12×76="
Replace-Text "35×97=" "18×25="
Replace-Text "33×66=" "84×98="
Replace-Text "19×78=" "16×54="
Replace-Text "58×46=" "59×37="
Replace-Text "93×69=" "39×31="
Replace-Text "23×88=" "16×70="
Replace-Text "52×19=" "63×61="
Replace-Text "50×44=" "75×59="
Replace-Text "97×59=" "46×68="
Replace-Text "97×81=" "44×87="
Replace-Text "27×94=" "61×60="
Replace-Text "14×55=" "77×86="
Replace-Text "42×46=" "95×14="
Replace-Text "72×20=" "35×31="
Replace-Text "80×34=" "84×37="
Replace-Text "82×32=" "46×14="
Replace-Text "67×55=" "71×30="
Replace-Text "26×41=" "99×70="
Replace-Text "62×98=" "48×76="
Replace-Text "42×75=" "61×91="
Replace-Text "17×40=" "68×19="
Replace-Text "50×80=" "21×19="
